# Update the workbook to reflect data through 2022-05-09
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab (and workbook sheet name) to reflect the new date
$ws.Name = "Through 2022-05-09"

# Update the "May (through 05-08)" label to "May (through 05-09)"
$ws.Range("A6").Value = "May (through 05-09)"

# Update the May row (row 6) values for years 2015-2022 (columns B-I)
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 11
$ws.Range("D6").Value = 18
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 18
$ws.Range("H6").Value = 35
$ws.Range("I6").Value = 31

# Update the Total row (row 7) values for years 2015-2022 (columns B-I)
$ws.Range("B7").Value = 97
$ws.Range("C7").Value = 173
$ws.Range("D7").Value = 271
$ws.Range("E7").Value = 254
$ws.Range("F7").Value = 165
$ws.Range("G7").Value = 280
$ws.Range("H7").Value = 558
$ws.Range("I7").Value = 583
